$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 18 with data for 2022-04-13 (date serial 44664)
$ws.Cells.Item(18, 1).Value = 44664
$ws.Cells.Item(18, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 326834
$ws.Cells.Item(18, 4).Value = 6338
$ws.Cells.Item(18, 5).Value = 31
$ws.Cells.Item(18, 6).Value = 1

# Update selection to match the author's final cursor position
$ws.Range("D18").Select()
